$d = $word.ActiveDocument

$replacements = @(
    @{Old = "27÷8="; New = "17÷7="},
    @{Old = "85÷9="; New = "87÷5="},
    @{Old = "48÷5="; New = "51÷4="},
    @{Old = "79÷3="; New = "59÷4="},
    @{Old = "50÷3="; New = "97÷5="},
    @{Old = "47÷8="; New = "62÷5="},
    @{Old = "23÷7="; New = "67÷5="},
    @{Old = "61÷4="; New = "44÷3="},
    @{Old = "42÷2="; New = "28÷2="},
    @{Old = "80÷8="; New = "29÷9="},
    @{Old = "97÷2="; New = "77÷7="},
    @{Old = "12÷2="; New = "69÷9="},
    @{Old = "99÷4="; New = "92÷4="},
    @{Old = "34÷8="; New = "45÷3="},
    @{Old = "37÷6="; New = "23÷6="},
    @{Old = "52÷7="; New = "27÷6="},
    @{Old = "73÷7="; New = "96÷9="},
    @{Old = "44÷5="; New = "39÷9="},
    @{Old = "66÷7="; New = "16÷5="},
    @{Old = "74÷3="; New = "23÷2="},
    @{Old = "75÷6="; New = "64÷9="},
    @{Old = "72÷2="; New = "35÷7="},
    @{Old = "31÷4="; New = "52÷5="},
    @{Old = "86÷7="; New = "38÷8="},
    @{Old = "26÷3="; New = "10÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
